$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.05405405405405406
$ws.Range("D3").Value = 0.3063063063063063
$ws.Range("E3").Value = 0.7207207207207207
$ws.Range("F3").Value = 0.9279279279279279
$ws.Range("H3").Value = 0.1378979670118911
$ws.Range("I3").Value = 0.3793458361359596
$ws.Range("J3").Value = -0.04504504504504504
$ws.Range("K3").Value = 1991.801801801802

$ws.Range("Q3").Value = 600
$ws.Range("R3").Value = 981
$ws.Range("S3").Value = 1656
$ws.Range("T3").Value = 2769
$ws.Range("U3").Value = 3541
$ws.Range("V3").Value = 4503
$ws.Range("W3").Value = 4122
$ws.Range("X3").Value = 3447
$ws.Range("Y3").Value = 2334
$ws.Range("Z3").Value = 1562

$ws.Range("AF3").Value = 0.882422
$ws.Range("AG3").Value = 0.80776
$ws.Range("AH3").Value = 0.675485
$ws.Range("AI3").Value = 0.457378
$ws.Range("AJ3").Value = 0.306094
